# Update numeric values in column F across the 4 worksheets, as described
# by the source diff (F-column "views/likes"-type counters incrementing).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1270
$ws1.Range("F4").Value  = 12937
$ws1.Range("F5").Value  = 737
$ws1.Range("F6").Value  = 61
$ws1.Range("F8").Value  = 60
$ws1.Range("F10").Value = 1878
$ws1.Range("F13").Value = 522
$ws1.Range("F15").Value = 125
$ws1.Range("F16").Value = 350
$ws1.Range("F18").Value = 300
$ws1.Range("F19").Value = 135
$ws1.Range("F20").Value = 131
$ws1.Range("F21").Value = 28
$ws1.Range("F22").Value = 221
$ws1.Range("F23").Value = 258
$ws1.Range("F24").Value = 1306
$ws1.Range("F25").Value = 341
$ws1.Range("F26").Value = 68

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 285
$ws2.Range("F5").Value  = 4453
$ws2.Range("F6").Value  = 164
$ws2.Range("F8").Value  = 16
$ws2.Range("F9").Value  = 74
$ws2.Range("F10").Value = 74
$ws2.Range("F11").Value = 360
$ws2.Range("F14").Value = 5

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 878
$ws3.Range("F3").Value = 4262

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 878
$ws4.Range("F6").Value  = 1271
$ws4.Range("F7").Value  = 12937
$ws4.Range("F8").Value  = 285
$ws4.Range("F9").Value  = 737
$ws4.Range("F10").Value = 4262
$ws4.Range("F11").Value = 61
$ws4.Range("F13").Value = 60
$ws4.Range("F15").Value = 1878
$ws4.Range("F18").Value = 522
$ws4.Range("F19").Value = 4453
$ws4.Range("F21").Value = 164
$ws4.Range("F22").Value = 164
$ws4.Range("F24").Value = 125
$ws4.Range("F25").Value = 16
$ws4.Range("F26").Value = 74
$ws4.Range("F27").Value = 74
$ws4.Range("F28").Value = 360
$ws4.Range("F29").Value = 350
$ws4.Range("F32").Value = 300
$ws4.Range("F33").Value = 135
$ws4.Range("F34").Value = 131
$ws4.Range("F35").Value = 28
$ws4.Range("F37").Value = 221
$ws4.Range("F38").Value = 5
$ws4.Range("F40").Value = 258
$ws4.Range("F41").Value = 1306
$ws4.Range("F43").Value = 341
$ws4.Range("F44").Value = 68
